$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 and 5 (old rows for MuSCs -> ECs and MuSCs -> MuSCs), leaving only
# the header row and the two remaining data rows.
$ws.Rows("4:5").Delete()

# Row 2: Sending cluster "ECs" -> "MuSCs" in Target cluster (D2), and update the
# downstream computed metrics to the new TPM-based values.
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 0.002166666666666667
$ws.Range("H2").Value = 0.0065
$ws.Range("I2").Value = 0.004890446475191893
$ws.Range("J2").Value = 0.004890446475191893
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.009620666666666666
$ws.Range("N2").Value = 0.028862
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.00002084477777777778
$ws.Range("R2").Value = 0.000187603
$ws.Range("S2").Value = 0.004890446475191893
$ws.Range("T2").Value = 0.004890446475191893

# Row 3: Sending cluster "ECs" -> "MuSCs" in A3, keep the rest as formerly row 3
# but update the metrics to new TPM-based values.
$ws.Range("A3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.440874
$ws.Range("H3").Value = 1.322622
$ws.Range("I3").Value = 0.9951095535248081
$ws.Range("J3").Value = 0.9951095535248081
$ws.Range("M3").Value = 0.009620666666666666
$ws.Range("N3").Value = 0.028862
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.004241501795999999
$ws.Range("R3").Value = 0.038173516164
$ws.Range("S3").Value = 0.9951095535248081
$ws.Range("T3").Value = 0.9951095535248081
